$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Coin name (B) and link (C) swaps/replacements ---
$textValues = @{
    "B35" = "HuobiToken"
    "B36" = "ImmutableX"
    "B47" = "EnergySwap"
    "B48" = "Cronos"
    "B51" = "EOS"
    "C35" = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
    "C36" = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
    "C47" = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
    "C48" = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
    "C51" = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
}
foreach ($ref in $textValues.Keys) {
    $ws.Range($ref).Value = $textValues[$ref]
}

# --- Volume / 1h change (E) — percentage text, safe from numeric coercion ---
$volumeValues = @{
    "E2" = "  +1.36%  "
    "E3" = "  +0.74%  "
    "E4" = "  -0.49%  "
    "E6" = "  -0.34%  "
    "E7" = "  +0.89%  "
    "E8" = "  +0.82%  "
    "E9" = "  +1.17%  "
    "E10" = "  -2.22%  "
    "E11" = "  +0.92%  "
    "E12" = "  +0.83%  "
    "E13" = "  +1.83%  "
    "E14" = "  +0.71%  "
    "E15" = "  -0.03%  "
    "E16" = "  -0.50%  "
    "E17" = "  -0.37%  "
    "E18" = "  +1.23%  "
    "E19" = "  -0.28%  "
    "E20" = "  +1.19%  "
    "E21" = "  +1.30%  "
    "E22" = "  +0.81%  "
    "E23" = "  -0.78%  "
    "E24" = "  +0.43%  "
    "E25" = "  -0.87%  "
    "E26" = "  -4.40%  "
    "E27" = "  +3.11%  "
    "E28" = "  -0.05%  "
    "E29" = "  +0.91%  "
    "E30" = "  -0.51%  "
    "E31" = "  -0.33%  "
    "E32" = "  +0.11%  "
    "E33" = "  -0.16%  "
    "E34" = "  +0.41%  "
    "E35" = "  +0.19%  "
    "E36" = "  +1.34%  "
    "E37" = "  +0.66%  "
    "E38" = "  -1.12%  "
    "E39" = "  -0.17%  "
    "E40" = "  -0.89%  "
    "E41" = "  +0.43%  "
    "E42" = "  -0.30%  "
    "E43" = "  +0.03%  "
    "E44" = "  -0.32%  "
    "E45" = "  -0.34%  "
    "E46" = "  -1.67%  "
    "E47" = "  -0.98%  "
    "E48" = "  +0.21%  "
    "E49" = "  -3.66%  "
    "E50" = "  -0.59%  "
    "E51" = "  +2.85%  "
}
foreach ($ref in $volumeValues.Keys) {
    $ws.Range($ref).Value = $volumeValues[$ref]
}

# --- Price (D) — force text format first so numeric-looking strings
#     (e.g. "0.9990", "26.737.41") are preserved exactly, matching the source feed
$priceValues = @{
    "D2" = "26.737.41"
    "D3" = "1.857.19"
    "D5" = "265.78"
    "D6" = "0.9990"
    "D7" = "0.5246"
    "D8" = "0.3292"
    "D9" = "0.06818"
    "D10" = "18.94"
    "D11" = "0.7793"
    "D12" = "0.07737"
    "D13" = "1.854.97"
    "D14" = "88.81"
    "D15" = "5.047"
    "D16" = "0.9992"
    "D17" = "14.08"
    "D19" = "0.9993"
    "D20" = "26.759.53"
    "D21" = "2.095.60"
    "D23" = "9.555"
    "D24" = "6.024"
    "D25" = "143.90"
    "D26" = "2.223"
    "D27" = "1.689"
    "D28" = "17.04"
    "D29" = "112.55"
    "D30" = "4.208"
    "D31" = "4.173"
    "D32" = "0.08766"
    "D33" = "0.04846"
    "D34" = "1.143"
    "D35" = "2.871"
    "D36" = "0.7174"
    "D37" = "3.122"
    "D39" = "2.214"
    "D40" = "0.4897"
    "D41" = "113.23"
    "D42" = "0.9028"
    "D43" = "6.097"
    "D44" = "0.9988"
    "D45" = "7.759"
    "D46" = "0.4207"
    "D47" = "9.143"
    "D48" = "0.05931"
    "D50" = "35.11"
    "D51" = "0.8873"
}
foreach ($ref in $priceValues.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $priceValues[$ref]
}
